$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values for rows 2-6 (historical data corrected)
# Row 2
$ws.Range("D2").Value = 1868
$ws.Range("E2").Value = 120
$ws.Range("F2").Value = 120
$ws.Range("G2").Value = 112
$ws.Range("H2").Value = 89
$ws.Range("I2").Value = 77
$ws.Range("J2").Value = 11
$ws.Range("K2").Value = 2507
$ws.Range("L2").Value = 707
$ws.Range("M2").Value = 1800
$ws.Range("N2").Value = 1686
$ws.Range("O2").Value = 114
$ws.Range("P2").Value = 70
$ws.Range("Q2").Value = 209
$ws.Range("R2").Value = -102
$ws.Range("S2").Value = -116
$ws.Range("T2").Value = 116
$ws.Range("U2").Value = 93
$ws.Range("V2").Value = 490
$ws.Range("W2").Value = 6.44
$ws.Range("X2").Value = 4.75
$ws.Range("Y2").Value = 4.7
$ws.Range("Z2").Value = 3.52
$ws.Range("AA2").Value = 39.28
$ws.Range("AB2").Value = 2309.55
$ws.Range("AC2").Value = 553
$ws.Range("AD2").Value = 9.03
$ws.Range("AE2").Value = 12042
$ws.Range("AF2").Value = 0.42
$ws.Range("AG2").Value = 50
$ws.Range("AH2").Value = 1
$ws.Range("AI2").Value = 9.03
$ws.Range("AJ2").Value = 14000000

# Row 3
$ws.Range("D3").Value = 1729
$ws.Range("E3").Value = 85
$ws.Range("F3").Value = 85
$ws.Range("G3").Value = 71
$ws.Range("H3").Value = 54
$ws.Range("I3").Value = 43
$ws.Range("J3").Value = 10
$ws.Range("K3").Value = 2941
$ws.Range("L3").Value = 1098
$ws.Range("M3").Value = 1844
$ws.Range("N3").Value = 1722
$ws.Range("O3").Value = 122
$ws.Range("P3").Value = 70
$ws.Range("Q3").Value = 53
$ws.Range("R3").Value = -123
$ws.Range("S3").Value = 146
$ws.Range("T3").Value = 130
$ws.Range("U3").Value = -77
$ws.Range("V3").Value = 649
$ws.Range("W3").Value = 4.9
$ws.Range("X3").Value = 3.11
$ws.Range("Y3").Value = 2.54
$ws.Range("Z3").Value = 1.98
$ws.Range("AA3").Value = 59.52
$ws.Range("AB3").Value = 2359.66
$ws.Range("AC3").Value = 309
$ws.Range("AD3").Value = 17.55
$ws.Range("AE3").Value = 12302
$ws.Range("AF3").Value = 0.44
$ws.Range("AG3").Value = 25
$ws.Range("AH3").Value = 0.46
$ws.Range("AI3").Value = 8.08
$ws.Range("AJ3").Value = 14000000

# Row 4
$ws.Range("D4").Value = 2253
$ws.Range("E4").Value = 202
$ws.Range("F4").Value = 202
$ws.Range("G4").Value = 323
$ws.Range("H4").Value = 249
$ws.Range("I4").Value = 207
$ws.Range("J4").Value = 42
$ws.Range("K4").Value = 3719
$ws.Range("L4").Value = 1483
$ws.Range("M4").Value = 2236
$ws.Range("N4").Value = 1981
$ws.Range("O4").Value = 255
$ws.Range("P4").Value = 70
$ws.Range("Q4").Value = 426
$ws.Range("R4").Value = -73
$ws.Range("S4").Value = -394
$ws.Range("T4").Value = 179
$ws.Range("U4").Value = 247
$ws.Range("V4").Value = 1004
$ws.Range("W4").Value = 8.96
$ws.Range("X4").Value = 11.05
$ws.Range("Y4").Value = 11.19
$ws.Range("Z4").Value = 7.48
$ws.Range("AA4").Value = 66.33
$ws.Range("AB4").Value = 2635.02
$ws.Range("AC4").Value = 1480
$ws.Range("AD4").Value = 5.97
$ws.Range("AE4").Value = 14151
$ws.Range("AF4").Value = 0.62
$ws.Range("AG4").Value = 75
$ws.Range("AH4").Value = 0.85
$ws.Range("AI4").Value = 5.07
$ws.Range("AJ4").Value = 14000000

# Row 5
$ws.Range("D5").Value = 2247
$ws.Range("E5").Value = 275
$ws.Range("F5").Value = 275
$ws.Range("G5").Value = 292
$ws.Range("H5").Value = 227
$ws.Range("I5").Value = 201
$ws.Range("J5").Value = 26
$ws.Range("K5").Value = 3636
$ws.Range("L5").Value = 1163
$ws.Range("M5").Value = 2474
$ws.Range("N5").Value = 2174
$ws.Range("O5").Value = 300
$ws.Range("P5").Value = 70
$ws.Range("Q5").Value = 310
$ws.Range("R5").Value = -115
$ws.Range("S5").Value = -236
$ws.Range("T5").Value = 218
$ws.Range("U5").Value = 91
$ws.Range("V5").Value = 739
$ws.Range("W5").Value = 12.24
$ws.Range("X5").Value = 10.08
$ws.Range("Y5").Value = 9.66
$ws.Range("Z5").Value = 6.16
$ws.Range("AA5").Value = 47.01
$ws.Range("AB5").Value = 2914.21
$ws.Range("AC5").Value = 1434
$ws.Range("AD5").Value = 6.28
$ws.Range("AE5").Value = 15530
$ws.Range("AF5").Value = 0.58
$ws.Range("AG5").Value = 100
$ws.Range("AH5").Value = 1.11
$ws.Range("AI5").Value = 6.97
$ws.Range("AJ5").Value = 14000000

# Row 6
$ws.Range("D6").Value = 2448
$ws.Range("E6").Value = 240
$ws.Range("F6").Value = 240
$ws.Range("G6").Value = 238
$ws.Range("H6").Value = 173
$ws.Range("I6").Value = 145
$ws.Range("K6").Value = 4062
$ws.Range("L6").Value = 1379
$ws.Range("M6").Value = 2683
$ws.Range("N6").Value = 2313
$ws.Range("P6").Value = 70
$ws.Range("Q6").Value = 101
$ws.Range("R6").Value = -139
$ws.Range("S6").Value = 133
$ws.Range("T6").Value = 92
$ws.Range("U6").Value = 9
$ws.Range("V6").Value = 921
$ws.Range("W6").Value = 9.81
$ws.Range("X6").Value = 7.06
$ws.Range("Y6").Value = 6.48
$ws.Range("Z6").Value = 4.49
$ws.Range("AA6").Value = 51.37
$ws.Range("AB6").Value = 3099.91
$ws.Range("AC6").Value = 1038
$ws.Range("AD6").Value = 7.56
$ws.Range("AE6").Value = 16521
$ws.Range("AF6").Value = 0.47
$ws.Range("AG6").Value = 50
$ws.Range("AH6").Value = 0.64
$ws.Range("AI6").Value = 4.82
$ws.Range("AJ6").Value = 14000000

# Rows 7-9: clear forecast columns except A, B, C (data no longer available)
$ws.Range("D7:E9").ClearContents()
$ws.Range("G7:I9").ClearContents()
$ws.Range("K7:N9").ClearContents()
$ws.Range("P7:U9").ClearContents()
$ws.Range("W7:AA9").ClearContents()
$ws.Range("AC7:AI9").ClearContents()
